# Update the "Model Comparison" table on slide 17 (RMSE values that were
# re-measured / corrected by the author).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$tbl = $s.Shapes.Item(2).Table

# Multinomial Logistic / LDA features : 0.702 -> 0.742
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "0.742"

# Ranger / Words Sparse Matrix + Length : 0.846 -> 0.793
$tbl.Cell(8, 3).Shape.TextFrame.TextRange.Text = "0.793"

# Ranger / LDA features : 0.886 -> 0.943
$tbl.Cell(9, 3).Shape.TextFrame.TextRange.Text = "0.943"

# Linear Regression / Words Sparse Matrix + Length : 0.855 -> 0.846
$tbl.Cell(10, 3).Shape.TextFrame.TextRange.Text = "0.846"

# Linear Regression / LDA features : 0.839 -> 0.965
$tbl.Cell(11, 3).Shape.TextFrame.TextRange.Text = "0.965"
